$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range so we know how many rows contain data
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Columns E (5) and F (6) need to be swapped for every row, including the
# header row: E holds "codeforiati:category-name" / category-name values and
# F holds "codeforiati:group-code" / group-code values before the edit; after
# the edit their contents are swapped (column E becomes group-code, column F
# becomes category-name).
$colE = 5
$colF = 6

for ($r = 1; $r -le $lastRow; $r++) {
    $cellE = $ws.Cells.Item($r, $colE)
    $cellF = $ws.Cells.Item($r, $colF)

    $valE = $cellE.Value2
    $valF = $cellF.Value2

    # Remember the cells' original style so it can be restored afterwards.
    $styleE = $cellE.Style
    $styleF = $cellF.Style

    # Every value in this sheet is stored as text (even the numeric-looking
    # codes such as "110"); temporarily force the number format to Text
    # before writing so Excel does not silently convert them back into
    # numbers, then restore the original style so no visible formatting
    # change is left behind.
    $cellE.NumberFormat = "@"
    $cellF.NumberFormat = "@"

    $cellE.Value = $valF
    $cellF.Value = $valE

    $cellE.Style = $styleE
    $cellF.Style = $styleF
}
